$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Finish the "units" table (rows 13-21) ---

# Row 13: add the column headers (Volume/Price/Weight/Surface) matching the
# style already used for A13 ("Ratio").
$ws.Range("A1").Copy()
$ws.Range("B13:E13").PasteSpecial(-4122)

$ws.Range("B13").Value = "Volume"
$ws.Range("C13").Value = "Price"
$ws.Range("D13").Value = "Weight"
$ws.Range("E13").Value = "Surface"

# Rows 14-21: replace the "X" placeholders with the actual operand-order
# numbers (1 / 2).
$ws.Range("B14").Value = 1
$ws.Range("E14").Value = 2

$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 2

$ws.Range("B16").Value = 2
$ws.Range("C16").Value = 1

$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 2

$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 2

$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 2

$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 2

$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 2

# --- View state updates ---
$ws.Application.ActiveWindow.Zoom = 225
$ws.Range("E21").Select()
